# "Refined metadata to be additional tab"
# Adds a new "metadata" worksheet (after "data") carrying the panel query
# metadata that used to live elsewhere, and refreshes the "time_taken"
# column on the "data" sheet to the new query timestamp.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the time_taken timestamps on the "data" sheet ---------------
$newTimes = @(
  "2021-10-05 14:20:32.894592",
  "2021-10-05 14:20:32.894600",
  "2021-10-05 14:20:32.894603",
  "2021-10-05 14:20:32.894606",
  "2021-10-05 14:20:32.894609",
  "2021-10-05 14:20:32.894611",
  "2021-10-05 14:20:32.894614",
  "2021-10-05 14:20:32.894616",
  "2021-10-05 14:20:32.894619",
  "2021-10-05 14:20:32.894621",
  "2021-10-05 14:20:32.894624",
  "2021-10-05 14:20:32.894626",
  "2021-10-05 14:20:32.894629",
  "2021-10-05 14:20:32.894631",
  "2021-10-05 14:20:32.894633",
  "2021-10-05 14:20:32.894636",
  "2021-10-05 14:20:32.894639",
  "2021-10-05 14:20:32.894641",
  "2021-10-05 14:20:32.894643",
  "2021-10-05 14:20:32.894646",
  "2021-10-05 14:20:32.894648",
  "2021-10-05 14:20:32.894651",
  "2021-10-05 14:20:32.894653",
  "2021-10-05 14:20:32.894655",
  "2021-10-05 14:20:32.894658",
  "2021-10-05 14:20:32.894661",
  "2021-10-05 14:20:32.894663",
  "2021-10-05 14:20:32.894666",
  "2021-10-05 14:20:32.894668",
  "2021-10-05 14:20:32.894670"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
  $row = $i + 2
  $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- Add the new "metadata" worksheet, placed after "data" --------------
$meta = $wb.Worksheets.Add($null, $dataSheet)
$meta.Name = "metadata"

# Re-use the "data" sheet's header style (bold, bordered, centered) for the
# new sheet's header row, and its row-2 "index" style for A2, by copying the
# formatted cells over before writing the new sheet's own values/text.
$dataSheet.Range("B1:F1").Copy($meta.Range("B1:G1"))
$dataSheet.Range("A2").Copy($meta.Range("A2"))

# Header row
$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

# Data row
$meta.Cells.Item(2, 1).Value = 0
$meta.Cells.Item(2, 2).Value = "GI tract tumours"
$meta.Cells.Item(2, 3).Value = 254

# "1.18" is a version string, not a number -- force text storage (so it
# doesn't get parsed/rounded as a float), then drop the number-format
# override again so the cell is left with the default (unstyled) look.
$meta.Range("D2").NumberFormat = "@"
$meta.Cells.Item(2, 4).Value = "1.18"
$meta.Range("D2").ClearFormats()

$meta.Cells.Item(2, 5).Value = "2019-08-05T14:17:21.117330Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:20:32.891076"
$meta.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/254/?format=json"

$dataSheet.Activate()
